# "Up to date with run 21"
#
# - Switch the active sheet from "FT-ELCP" to "Demand", updating the
#   Demand sheet's current selection to D13.
# - On the "Demand" sheet, change the annual growth-rate input (I4) from
#   10% to 1%. The projected demand figures in G15:G34 are driven by
#   formulas that reference I4, so they recalculate automatically.

$wb = $excel.ActiveWorkbook

$wsDemand = $wb.Worksheets.Item("Demand")

# Make "Demand" the active/selected sheet (this also un-selects the
# previously active "FT-ELCP" sheet) and set its new selection.
$wsDemand.Activate()
$wsDemand.Range("D13").Select()

# Update the growth rate assumption; dependent formulas recalc automatically.
$wsDemand.Range("I4").Formula = "=1%"
